$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304" ---
$oldHeaders = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $oldHeaders[$i] + "_FV2210"
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = $oldHeaders[$i] + "_FV2304"
}

# --- Turn the used range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- Freeze the header row (split/freeze after row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
